# "List Kebutuhan.xlsx" - Sheet1 update:
#   - Section 6/7 rewritten from admin-centric wording to the generic wording
#   - F4.3 wording tweak ("transaksi reservasi" -> "transaksi peminjaman ruangan")
#   - Column C narrowed (68.78 -> 58 chars) now that the text is shorter
#   - Rows 14/15 shrink back to the default row height once re-autofit
#   - Selection/viewport nudged down to the edited area (C12)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# F 4.3 wording tweak
$ws.Range("C11").Value = "F 4.3 Sistem memungkinkan pengguna yang telah memesan ruangan untuk melakukan booking pada transaksi peminjaman ruangan"

# Section 6 ("Melayani pemesanan ruangan") - simplified, no longer admin-specific
$ws.Range("B14").Value = "6. Melayani pemesanan ruangan"
$ws.Range("C14").Value = "F 6.1 Menambah informasi ruangan"
$ws.Range("C15").Value = "F 6.2 Melihat ketersediaan ruangan"
$ws.Range("C16").Value = "F 6.3 Melihat data diri peminjam ruangan"

# Section 7 ("Melakukan komunikasi dengan peminjam ruangan") - simplified
$ws.Range("B17").Value = "7. Melakukan komunikasi dengan peminjam ruangan"
$ws.Range("C17").Value = "F 7.1 Melakukan chat dengan peminjam ruangan"

# Column C can be narrower now that the Section 6/7 text is shorter
$ws.Columns.Item(3).ColumnWidth = 57.17

# Let rows 14-15 re-fit to the (now shorter) text
$ws.Range("A14:D15").EntireRow.AutoFit()

# Move the selection/viewport down near the edited rows
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("C12").Select()
